$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: A2 0, B2 281
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 281

# Update row 3: A3 1, B3 stays 175
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 175

# Remove row 4 (A4=2, B4=56) entirely
$ws.Range("A4:B4").Delete()
